$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "in"/"out" header columns: B1 becomes "out", C1 becomes "in"
$ws.Range("B1").Value = "out"
$ws.Range("C1").Value = "in"

# Update existing data rows (2-7) with new in/out numbers
$ws.Range("B2").Value = 288
$ws.Range("C2").Value = 445

$ws.Range("B3").Value = 370
$ws.Range("C3").Value = 402

$ws.Range("B4").Value = 198
$ws.Range("C4").Value = 256

$ws.Range("B6").Value = 117
$ws.Range("C6").Value = 167

$ws.Range("B7").Value = 242
$ws.Range("C7").Value = 353

# Add two new rows for additional leagues
$ws.Range("A8").Value = "Eredivisie"
$ws.Range("B8").Value = 162
$ws.Range("C8").Formula = "=400-B8"

$ws.Range("A9").Value = "Champioship"
$ws.Range("B9").Value = 2018
$ws.Range("C9").Value = 1790

# Match the selection recorded in the saved file
$ws.Range("C2").Select()
